# Update "Generate Report for Handback" timestamps for the
# 6419ced1-b836-49b0-bd9d-e5c48dc6a783 entry (row 4) on the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 -> D4 (Correspond Handoff Datetime), G4 (Correspond Handback DateTime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-19 07:07:13"
$wsZhCn.Range("G4").Value = "2016-02-19 07:07:56"

# de-de sheet: row 4 -> D4 (Correspond Handoff Datetime), G4 (Correspond Handback DateTime)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-19 07:07:23"
$wsDeDe.Range("G4").Value = "2016-02-19 07:08:12"
